$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to text
# (matching the source inlineStr cell type) by pre-setting a text number format.
$textCells = @('D5','D6','D8','D11','D18','D20','D23','D25','D27','D28','D34','D36','D39','D42','D45','D48','D50','D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '27.460.42'
$ws.Range('E2').Value = '  -0.86%  '
$ws.Range('D3').Value = '1.567.58'
$ws.Range('E3').Value = '  -1.10%  '
$ws.Range('E4').Value = '  -0.22%  '
$ws.Range('D5').Value = '208.49'
$ws.Range('E5').Value = '  +0.75%  '
$ws.Range('D6').Value = '0.501'
$ws.Range('E6').Value = '  -0.99%  '
$ws.Range('E7').Value = '  -0.21%  '
$ws.Range('D8').Value = '22.16'
$ws.Range('E8').Value = '  -0.42%  '
$ws.Range('E9').Value = '  -1.30%  '
$ws.Range('E10').Value = '  +0.31%  '
$ws.Range('D11').Value = '0.0865'
$ws.Range('E11').Value = '  -0.34%  '
$ws.Range('D12').Value = '1.790.86'
$ws.Range('E12').Value = '  -1.03%  '
$ws.Range('D13').Value = '1.557.24'
$ws.Range('E13').Value = '  -2.85%  '
$ws.Range('E14').Value = '  -0.86%  '
$ws.Range('E15').Value = '  -2.47%  '
$ws.Range('E16').Value = '  +0.71%  '
$ws.Range('D17').Value = '27.464.97'
$ws.Range('E17').Value = '  -0.70%  '
$ws.Range('D18').Value = '214.08'
$ws.Range('D19').Value = '0.0₃0691'
$ws.Range('E19').Value = '  -0.30%  '
$ws.Range('D20').Value = '7.27'
$ws.Range('E20').Value = '  -0.85%  '
$ws.Range('E21').Value = '  -0.18%  '
$ws.Range('E22').Value = '  -0.69%  '
$ws.Range('D23').Value = '9.56'
$ws.Range('E23').Value = '  +0.66%  '
$ws.Range('E24').Value = '  +2.01%  '
$ws.Range('D25').Value = '152.86'
$ws.Range('E25').Value = '  -1.17%  '
$ws.Range('E26').Value = '  -0.21%  '
$ws.Range('D27').Value = '6.71'
$ws.Range('E27').Value = '  -1.79%  '
$ws.Range('D28').Value = '15.00'
$ws.Range('E28').Value = '  -0.79%  '
$ws.Range('E29').Value = '  -1.76%  '
$ws.Range('E30').Value = '  -0.50%  '
$ws.Range('E31').Value = '  +1.05%  '
$ws.Range('E32').Value = '  -0.89%  '
$ws.Range('D33').Value = '1.378.30'
$ws.Range('E33').Value = '  -0.26%  '
$ws.Range('D34').Value = '2.98'
$ws.Range('E34').Value = '  +1.91%  '
$ws.Range('E35').Value = '  +1.35%  '
$ws.Range('D36').Value = '0.953'
$ws.Range('E36').Value = '  -1.63%  '
$ws.Range('E37').Value = '  -0.78%  '
$ws.Range('E38').Value = '  +1.39%  '
$ws.Range('D39').Value = '0.543'
$ws.Range('E39').Value = '  +1.18%  '
$ws.Range('E40').Value = '  +0.99%  '
$ws.Range('E41').Value = '  -0.18%  '
$ws.Range('D42').Value = '0.980'
$ws.Range('E42').Value = '  +0.61%  '
$ws.Range('E43').Value = '  +3.13%  '
$ws.Range('D45').Value = '2.17'
$ws.Range('E45').Value = '  -0.19%  '
$ws.Range('E46').Value = '  +1.14%  '
$ws.Range('D47').Value = '1.703.34'
$ws.Range('D48').Value = '85.22'
$ws.Range('E48').Value = '  -3.23%  '
$ws.Range('E49').Value = '  +0.39%  '
$ws.Range('D50').Value = '0.0959'
$ws.Range('E50').Value = '  -1.28%  '
$ws.Range('D51').Value = '0.0496'
$ws.Range('E51').Value = '  -0.74%  '

# Restore default cell style on the forced-text cells so formatting matches the original.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
